$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.839.05'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.812.72'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4322'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3713'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07289'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8683'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.11%  '
$ws.Range("D12").Value = '1.996.93'
$ws.Range("E12").Value = '  +8.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.639'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.365'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06945'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008929'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").Value = '26.857.74'
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.211'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '2.185.69'
$ws.Range("E24").Value = '  +6.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.869'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.243'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.893'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08941'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7585'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.177'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.448'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.803'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.006'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.130'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05241'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01925'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5090'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1651'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.666'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.556'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.334'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.656'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4585'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06294'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.817'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.49%  '
